$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# Row 39: new part 1031 - Power and Signal Board Case Bottom
$ws.Range("A39").Value = 1031
$ws.Range("C39").Value = "3D Printed"
$ws.Range("E39").Value = "Power and Signal Board Case Bottom"

# Row 40: new part 1032 - Power and signal board case top
$ws.Range("A40").Value = 1032
$ws.Range("C40").Value = "3D Printed"
$ws.Range("E40").Value = "Power and signal board case top"

# Row 41: stray cell with a single "]" value (matches source edit)
$ws.Range("A41").Value = "]"

$ws.Range("A41").Select()
